$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.305.07"
$ws.Range("E2").Value = "  -3.44%  "

$ws.Range("D3").Value = "3.140.36"
$ws.Range("E3").Value = "  -5.04%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.44"
$ws.Range("E5").Value = "  -6.10%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.66"
$ws.Range("E6").Value = "  -5.43%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.140.63"
$ws.Range("E8").Value = "  -4.99%  "

$ws.Range("E9").Value = "  -6.58%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.21"
$ws.Range("E10").Value = "  -8.26%  "

$ws.Range("E11").Value = "  -9.18%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("E12").Value = "  -7.41%  "

$ws.Range("D13").Value = "3.677.27"
$ws.Range("E13").Value = "  -5.03%  "

$ws.Range("E14").Value = "  -1.82%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.47"
$ws.Range("E15").Value = "  -5.52%  "

$ws.Range("D16").Value = "3.138.28"
$ws.Range("E16").Value = "  -5.05%  "

$ws.Range("D17").Value = "58.261.57"
$ws.Range("E17").Value = "  -3.52%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000152"
$ws.Range("E18").Value = "  -8.51%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.78"
$ws.Range("E19").Value = "  -5.53%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.04"
$ws.Range("E20").Value = "  -7.22%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.92"
$ws.Range("E21").Value = "  -8.63%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "342.81"
$ws.Range("E22").Value = "  -8.72%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "67.76"
$ws.Range("E24").Value = "  -8.94%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.505"
$ws.Range("E25").Value = "  -5.92%  "

$ws.Range("D26").Value = "3.260.20"
$ws.Range("E26").Value = "  -5.12%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.167"
$ws.Range("E27").Value = "  -2.93%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "0.0₃0954"
$ws.Range("E28").Value = "  -7.19%  "

$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("E29").Value = "  -0.18%  "

$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.79"
$ws.Range("E30").Value = "  -5.71%  "

$ws.Range("B31").Value = "USDe"
$ws.Range("C31").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("E32").Value = "  -8.62%  "

$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.79"
$ws.Range("E33").Value = "  -10.32%  "

$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.41"
$ws.Range("E34").Value = "  -5.48%  "

$ws.Range("E35").Value = "  -3.25%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.80"
$ws.Range("E36").Value = "  -6.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "156.72"
$ws.Range("E37").Value = "  -6.12%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.21"
$ws.Range("E38").Value = "  -7.29%  "

$ws.Range("E39").Value = "  -11.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0687"
$ws.Range("E40").Value = "  -6.07%  "

$ws.Range("D41").Value = "3.171.13"
$ws.Range("E41").Value = "  -4.99%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "40.37"
$ws.Range("E42").Value = "  -3.73%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "24.09"
$ws.Range("E43").Value = "  -9.28%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.691"
$ws.Range("E44").Value = "  -8.02%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.08"
$ws.Range("E45").Value = "  -3.11%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.89"
$ws.Range("E46").Value = "  -6.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  +0.04%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.43"
$ws.Range("E48").Value = "  -9.37%  "

$ws.Range("B49").Value = "Maker"
$ws.Range("C49").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D49").Value = "2.260.31"
$ws.Range("E49").Value = "  -3.97%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.19"
$ws.Range("E50").Value = "  -3.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.61"
$ws.Range("E51").Value = "  -3.68%  "
